# "add Use Item bug" — NPC.xlsx / Property sheet:
#   - Row 49 (was "ShowName") becomes "EquipIDRef" (moved up, loses its old
#     custom row styling so it matches the plain style already used by the
#     row below it).
#   - Row 50 becomes "ShowName" (i.e. ShowName/EquipIDRef swap places).
#   - Two new rows are appended: "Icon" (51) and "ShowCard" (52).
#
# Columns (row 1 header): A=Id B=Type C=Public D=Private E=Save F=View
#                          G=Index H=SaveInterval I=RelationValue J=Desc

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 previously carried its own custom row/cell styling (s="7" row,
# s="6" cells) inherited from the old "ShowName" row. Drop that so the row
# matches the plain (unstyled) look the diff shows once it becomes
# "EquipIDRef" — same look row 50 already had.
$ws.Rows(49).ClearFormats()

# Row 49: EquipIDRef
$ws.Range("A49").Value = "EquipIDRef"
$ws.Range("B49").Value = "string"
$ws.Range("C49").Value = $false
$ws.Range("D49").Value = $false
$ws.Range("E49").Value = $false
$ws.Range("F49").Value = $true
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = "Friend"
$ws.Range("I49").NumberFormat = "@"
$ws.Range("J49").Value = "装备列表索引"

# Row 50: ShowName
$ws.Range("A50").Value = "ShowName"
$ws.Range("B50").Value = "string"
$ws.Range("C50").Value = $false
$ws.Range("D50").Value = $false
$ws.Range("E50").Value = $false
$ws.Range("F50").Value = $true
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = "Friend"
$ws.Range("I50").NumberFormat = "@"
$ws.Range("J50").Value = "显示名字"

# Row 51 (new): Icon
$ws.Range("A51").Value = "Icon"
$ws.Range("B51").Value = "string"
$ws.Range("C51").Value = $false
$ws.Range("D51").Value = $false
$ws.Range("E51").Value = $false
$ws.Range("F51").Value = $true
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = "Friend"
$ws.Range("I51").NumberFormat = "@"
$ws.Range("J51").Value = "图标"

# Row 52 (new): ShowCard
$ws.Range("A52").Value = "ShowCard"
$ws.Range("B52").Value = "string"
$ws.Range("C52").Value = $false
$ws.Range("D52").Value = $false
$ws.Range("E52").Value = $false
$ws.Range("F52").Value = $true
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = "Friend"
$ws.Range("I52").NumberFormat = "@"
$ws.Range("J52").Value = "卡牌背景"

# Match the author's final selection / scroll position as closely as the
# object model allows.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("A48").Select()
